$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text updates (Volume number, report week dates) ---
$ws.Range("A8").Value = "Volume 29   Number  40"
$ws.Range("C9").Value = "Report Covering the Week  10/3/2022  Through  10/9/2022"

# --- Row 14 ---
$ws.Range("L14").Value = -72.727272727272

# --- Row 15 ---
$ws.Range("D15").Value = 2
$ws.Range("E15").Value = -50
$ws.Range("F15").Value = 5
$ws.Range("G15").Value = 4
$ws.Range("H15").Value = 25
$ws.Range("I15").Value = 31
$ws.Range("J15").Value = 22
$ws.Range("K15").Value = 40.909090909090
$ws.Range("L15").Value = 34.782608695652
$ws.Range("M15").Value = 19.230769230769
$ws.Range("N15").Value = 24

# --- Row 16 ---
$ws.Range("D16").Value = 6
$ws.Range("E16").Value = 66.666666666666
$ws.Range("F16").Value = 31
$ws.Range("G16").Value = 22
$ws.Range("H16").Value = 40.909090909090
$ws.Range("I16").Value = 230
$ws.Range("J16").Value = 190
$ws.Range("K16").Value = 21.052631578947
$ws.Range("L16").Value = 50.326797385620
$ws.Range("M16").Value = -13.533834586466
$ws.Range("N16").Value = -77.648202137998

# --- Row 17 ---
$ws.Range("C17").Value = 9
$ws.Range("E17").Value = -10
$ws.Range("F17").Value = 27
$ws.Range("G17").Value = 41
$ws.Range("H17").Value = -34.146341463414
$ws.Range("I17").Value = 331
$ws.Range("J17").Value = 303
$ws.Range("K17").Value = 9.240924092409
$ws.Range("L17").Value = 35.102040816326
$ws.Range("M17").Value = 27.307692307692
$ws.Range("N17").Value = -1.780415430267

# --- Row 18 ---
$ws.Range("C18").Value = 4
$ws.Range("D18").Value = 3
$ws.Range("E18").Value = 33.333333333333
$ws.Range("F18").Value = 11
$ws.Range("G18").Value = 14
$ws.Range("H18").Value = -21.428571428571
$ws.Range("I18").Value = 118
$ws.Range("J18").Value = 99
$ws.Range("K18").Value = 19.191919191919
$ws.Range("L18").Value = -9.923664122137
$ws.Range("M18").Value = -48.017621145374
$ws.Range("N18").Value = -92.401802962009

# --- Row 19 ---
$ws.Range("C19").Value = 21
$ws.Range("E19").Value = 40
$ws.Range("F19").Value = 72
$ws.Range("H19").Value = 46.938775510204
$ws.Range("I19").Value = 755
$ws.Range("J19").Value = 409
$ws.Range("K19").Value = 84.596577017114
$ws.Range("L19").Value = 87.810945273631
$ws.Range("M19").Value = 98.162729658792
$ws.Range("N19").Value = -33.597185576077

# --- Row 20 ---
$ws.Range("C20").Value = 6
$ws.Range("E20").Value = 0
$ws.Range("F20").Value = 30
$ws.Range("G20").Value = 21
$ws.Range("H20").Value = 42.857142857142
$ws.Range("I20").Value = 248
$ws.Range("J20").Value = 162
$ws.Range("K20").Value = 53.086419753086
$ws.Range("L20").Value = 61.038961038961
$ws.Range("M20").Value = 41.714285714285
$ws.Range("N20").Value = -85.836664762992

# --- Row 21 ---
$ws.Range("C21").Value = 51
$ws.Range("D21").Value = 42
$ws.Range("E21").Value = 21.428571428571
$ws.Range("F21").Value = 176
$ws.Range("G21").Value = 151
$ws.Range("H21").Value = 16.556291390728
$ws.Range("I21").Value = 1716
$ws.Range("J21").Value = 1189
$ws.Range("K21").Value = 44.322960470984
$ws.Range("L21").Value = 53.351206434316
$ws.Range("M21").Value = 28.443113772455
$ws.Range("N21").Value = -70.651616213442

# --- Row 22 ---
$ws.Range("C22").Value = 2
$ws.Range("G22").Value = 3
$ws.Range("H22").Value = 133.333333333333
$ws.Range("I22").Value = 49
$ws.Range("K22").Value = 122.727272727273
$ws.Range("L22").Value = 226.666666666667
$ws.Range("M22").Value = 81.481481481481
$ws.Range("D22").Value = "'0"
$ws.Range("D23").Copy()
$ws.Range("D22").PasteSpecial(-4122)
$ws.Range("E22").Value = "'***.*"
$ws.Range("E23").Copy()
$ws.Range("E22").PasteSpecial(-4122)

# --- Row 24 ---
$ws.Range("C24").Value = 29
$ws.Range("D24").Value = 28
$ws.Range("E24").Value = 3.571428571428
$ws.Range("F24").Value = 122
$ws.Range("G24").Value = 109
$ws.Range("H24").Value = 11.926605504587
$ws.Range("I24").Value = 1436
$ws.Range("J24").Value = 1044
$ws.Range("K24").Value = 37.547892720306
$ws.Range("L24").Value = 34.835680751173
$ws.Range("M24").Value = 65.057471264367

# --- Row 25 ---
$ws.Range("C25").Value = 16
$ws.Range("D25").Value = 15
$ws.Range("E25").Value = 6.666666666666
$ws.Range("F25").Value = 63
$ws.Range("G25").Value = 67
$ws.Range("H25").Value = -5.970149253731
$ws.Range("I25").Value = 686
$ws.Range("J25").Value = 639
$ws.Range("K25").Value = 7.355242566510
$ws.Range("L25").Value = 16.865417376490
$ws.Range("M25").Value = -7.046070460704

# --- Row 26 ---
$ws.Range("F26").Value = 5
$ws.Range("G26").Value = 6
$ws.Range("H26").Value = -16.666666666666
$ws.Range("I26").Value = 40
$ws.Range("J26").Value = 40
$ws.Range("K26").Value = 0
$ws.Range("L26").Value = 11.111111111111

# --- Row 27 ---
$ws.Range("C27").Value = 3
$ws.Range("D27").Value = 2
$ws.Range("E27").Value = 50
$ws.Range("F27").Value = 10
$ws.Range("H27").Value = 150
$ws.Range("I27").Value = 80
$ws.Range("J27").Value = 81
$ws.Range("K27").Value = -1.234567901234
$ws.Range("L27").Value = 17.647058823529

# --- Row 28 ---
$ws.Range("L28").Value = 14.285714285714
$ws.Range("N28").Value = -83.673469387755

# --- Row 29 ---
$ws.Range("L29").Value = 16.666666666666
$ws.Range("N29").Value = -84.444444444444

# --- Row 30 ---
$ws.Range("G30").Value = 1
$ws.Range("H30").Value = 0
$ws.Range("C30").Value = "'0"
$ws.Range("C23").Copy()
$ws.Range("C30").PasteSpecial(-4122)
